$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bibi Cell Mundi)
$ws.Range("H2").Value = 16224.4
$ws.Range("I2").Value = 22772.55
$ws.Range("J2").Value = 2365
$ws.Range("AG2").Value = 86631.73

# Row 3 (Bibi Cell Vieiralves)
$ws.Range("H3").Value = 6192
$ws.Range("I3").Value = 9352
$ws.Range("J3").Value = 5889
$ws.Range("AG3").Value = 52385.2

# Row 4 (Bibi Cell Manauara)
$ws.Range("H4").Value = 3537
$ws.Range("I4").Value = 3125
$ws.Range("J4").Value = 4680.4
$ws.Range("K4").Value = 1488
$ws.Range("AG4").Value = 30839.15

# Row 5 (Bibi Cell Ponta Negra)
$ws.Range("H5").Value = 2554
$ws.Range("I5").Value = 2532.9
$ws.Range("J5").Value = 3252.89
$ws.Range("K5").Value = 1209.99
$ws.Range("AG5").Value = 25317.33

# Row 6 (total)
$ws.Range("H6").Value = 28507.4
$ws.Range("I6").Value = 37782.45
$ws.Range("J6").Value = 16187.29
$ws.Range("K6").Value = 2697.99
$ws.Range("AG6").Value = 195173.41
